$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New column BY: "taxa_barstacks" project-report property -------------
# Row 1 -> long description (header style copied from A1/B1, style index 1)
# Row 2 -> field name "taxa_barstacks" (style copied from A2/B2, style index 2)
# Rows 3-7 -> default value "phylum, class" for every demo sample (style copied
#             from B3, style index 4 - same as the other plain value columns)

# Shared-string table is appended to in reference order, so write the cells
# in the same order the strings are expected to appear (field name, then
# default value, then the long description).
$ws.Range("BY2").Value = "taxa_barstacks"
$ws.Range("BY3").Value = "phylum, class"
$ws.Range("BY4").Value = "phylum, class"
$ws.Range("BY5").Value = "phylum, class"
$ws.Range("BY6").Value = "phylum, class"
$ws.Range("BY7").Value = "phylum, class"
$ws.Range("BY1").Value = "This property determines how many taxonomical barstacks plots are included in the project report. Seperate each value with a comma. Possible values are: ""phylum, class, order, family, genus, species"""

# Match formatting to the equivalent cells used by every other column.
$ws.Range("A1").Copy()
$ws.Range("BY1").PasteSpecial(-4122)

$ws.Range("A2").Copy()
$ws.Range("BY2").PasteSpecial(-4122)

$ws.Range("B3").Copy()
$ws.Range("BY3:BY7").PasteSpecial(-4122)

# Give the new column the same kind of custom width the other description
# columns have (stored width 24.5 once Excel applies its internal padding).
$ws.Columns.Item(77).ColumnWidth = 23.666666666666668

# Reflect the cursor move that happened while adding this column.
$null = $ws.Range("BU17").Select()
